# KP-11725 D: Extension of questionnaire's translation files
# Adds a new "Variable" column (containing the question/option variable name,
# here "c1") right after the "Entity Id" column on both worksheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Translations")
$ws2 = $wb.Worksheets.Item("@@_question")

# --- Insert the new "Variable" column (becomes column B) on both sheets ----
$ws1.Columns.Item(2).Insert()
$ws2.Columns.Item(2).Insert()

# --- Header ------------------------------------------------------------
$ws1.Range("B1").Value = "Variable"
$ws2.Range("B1").Value = "Variable"

# --- Data rows: variable name "c1" for every data row -------------------
$ws1.Range("B2").Value = "c1"
$ws1.Range("B3").Value = "c1"
$ws1.Range("B4").Value = "c1"
$ws1.Range("B5").Value = "c1"

$ws2.Range("B2").Value = "c1"

# --- Approximate the column width of the freshly inserted column --------
$ws1.Columns.Item(2).ColumnWidth = 12.830729166666666
$ws2.Columns.Item(2).ColumnWidth = 6.830729166666667

# --- Restore view/selection state ---------------------------------------
# "@@_question" keeps a selection on the (now shifted) Variable column...
$ws2.Activate()
$ws2.Range("B2").Select() | Out-Null

# ...while "Translations" becomes the active/selected tab again.
$ws1.Activate()
$ws1.Range("B6").Select() | Out-Null
